$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '65.594.49'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '3.278.55'
$ws.Range("E3").Value = '  -1.13%  '
Set-TextValue "D4" '0.998'
$ws.Range("E4").Value = '  -0.15%  '
Set-TextValue "D5" '575.87'
$ws.Range("E5").Value = '  +2.91%  '
Set-TextValue "D6" '181.48'
$ws.Range("E6").Value = '  -3.63%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = '3.274.96'
$ws.Range("E8").Value = '  -0.99%  '
Set-TextValue "D9" '0.569'
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("E10").Value = '  -6.35%  '
$ws.Range("E11").Value = '  -3.18%  '
Set-TextValue "D12" '46.17'
$ws.Range("E12").Value = '  -3.22%  '
Set-TextValue "D13" '0.0000262'
$ws.Range("E13").Value = '  -3.71%  '
$ws.Range("D14").Value = '3.801.44'
$ws.Range("E14").Value = '  -1.13%  '
Set-TextValue "D15" '8.35'
$ws.Range("E15").Value = '  -3.36%  '
Set-TextValue "D16" '609.85'
$ws.Range("E16").Value = '  -3.57%  '
$ws.Range("D17").Value = '65.685.02'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D19" '17.61'
$ws.Range("E19").Value = '  -3.13%  '
$ws.Range("B20").Value = 'WrappedEther'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D20").Value = '3.278.56'
$ws.Range("E20").Value = '  -1.26%  '
Set-TextValue "D21" '10.84'
$ws.Range("E21").Value = '  -2.83%  '
Set-TextValue "D22" '0.883'
$ws.Range("E22").Value = '  -2.88%  '
Set-TextValue "D23" '18.25'
$ws.Range("E23").Value = '  +0.89%  '
Set-TextValue "D24" '4.91'
$ws.Range("E24").Value = '  -0.59%  '
Set-TextValue "D25" '97.89'
$ws.Range("E25").Value = '  -5.27%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  -0.73%  '
Set-TextValue "D28" '9.39'
$ws.Range("E28").Value = '  -1.95%  '
Set-TextValue "D29" '30.58'
$ws.Range("E29").Value = '  +0.88%  '
Set-TextValue "D30" '8.34'
$ws.Range("E30").Value = '  -4.02%  '
$ws.Range("E31").Value = '  +0.83%  '
Set-TextValue "D32" '3.72'
$ws.Range("E32").Value = '  -8.50%  '
Set-TextValue "D33" '544.97'
$ws.Range("E33").Value = '  -1.47%  '
Set-TextValue "D34" '10.79'
$ws.Range("E34").Value = '  -2.69%  '
$ws.Range("D35").Value = '3.793.12'
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("E37").Value = '  -0.05%  '
Set-TextValue "D38" '56.02'
$ws.Range("E38").Value = '  -2.92%  '
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D40" '32.38'
$ws.Range("E40").Value = '  -4.08%  '
$ws.Range("B41").Value = 'ApeXProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue "D41" '3.40'
$ws.Range("E41").Value = '  +4.42%  '
Set-TextValue "D42" '3.13'
$ws.Range("E42").Value = '  -5.15%  '
$ws.Range("D43").Value = '0.0₃0677'
$ws.Range("E43").Value = '  -7.84%  '
$ws.Range("E44").Value = '  -4.98%  '
Set-TextValue "D45" '0.328'
$ws.Range("E45").Value = '  -2.05%  '
Set-TextValue "D46" '0.0403'
$ws.Range("E46").Value = '  -4.06%  '
Set-TextValue "D47" '2.99'
$ws.Range("E47").Value = '  -6.95%  '
$ws.Range("E48").Value = '  +0.25%  '
Set-TextValue "D49" '0.126'
$ws.Range("E49").Value = '  -2.83%  '
$ws.Range("E50").Value = '  -4.67%  '
Set-TextValue "D51" '127.76'
$ws.Range("E51").Value = '  +4.64%  '
